# Fruta / hortaliza, semanal
# A new weekly price record was inserted as row 400 (shifting the existing
# rows 400-448 down to 401-449).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 400, pushing rows 400..448 down to 401..449.
$ws.Rows.Item(400).Insert()

# Populate the newly inserted row 400 with the new weekly record.
$ws.Cells.Item(400, 1).Value = 8
$ws.Cells.Item(400, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(400, 3).Value = "Coquimbo"
$ws.Cells.Item(400, 4).Value = 44918
$ws.Cells.Item(400, 5).Value = 4
$ws.Cells.Item(400, 6).Value = 100114013
$ws.Cells.Item(400, 7).Value = "Zanahoria"
$ws.Cells.Item(400, 8).Value = "Sin especificar"
$ws.Cells.Item(400, 9).Value = "Primera"
$ws.Cells.Item(400, 10).Value = 440
$ws.Cells.Item(400, 11).Value = 5000
$ws.Cells.Item(400, 12).Value = 6000
$ws.Cells.Item(400, 13).Value = 5500
$ws.Cells.Item(400, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(400, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(400, 16).Value = 275
$ws.Cells.Item(400, 17).Value = 20
$ws.Cells.Item(400, 18).Value = "Hortaliza"
